# Update cryptocurrency price/volume data per the scraper refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.598.97"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +6.46%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.733.94"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.26%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9959"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.18"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +5.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9937"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3715"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.58"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.89%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3381"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.204"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.06%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07528"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +6.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9941"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.428"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +5.87%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.49"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.018"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +5.68%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.726.73"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.69%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001094"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.80%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06679"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.87%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.59"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.97%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9938"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.78"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +6.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.179"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.17"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +4.34%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "26.532.42"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +6.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.462"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.526"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.426"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +15.98%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.00"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.09%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.77%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.919.86"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "131.66"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.113"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.147"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08593"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.65%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.21"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.702"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.454"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.07%  "

$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06359"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.93%  "

$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02354"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.58%  "

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2173"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.60%  "

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.675"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.53%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.242"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6264"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.93%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.54"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +13.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9941"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.901"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6074"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +7.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.56"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.064"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07340"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.97"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.80%  "
